$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.00029166883905418217
$ws.Range("A3").Value = 0.0001398509048158303
$ws.Range("H3").Value = 6.290065765380859
$ws.Range("A4").Value = 0.00011952534259762615
$ws.Range("H4").Value = 5.308176040649414
